# Fruta / hortaliza, semanal
# Insert two new weekly rows (268, 269) above the former row 268, pushing
# the existing rows 268-281 down to 270-283. Populate the two new rows
# with the latest week's data (fecha 44509, origen "Limache").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 268; existing rows 268..281 shift to 270..283.
$ws.Rows("268:269").Insert()

# --- New row 268 ---
$ws.Cells.Item(268, 1).Value2 = 11
$ws.Cells.Item(268, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(268, 3).Value2 = "Bíobío"
$ws.Cells.Item(268, 4).Value2 = 44509
$ws.Cells.Item(268, 5).Value2 = 8
$ws.Cells.Item(268, 6).Value2 = 100112020
$ws.Cells.Item(268, 7).Value2 = "Tomate"
$ws.Cells.Item(268, 8).Value2 = "Larga vida"
$ws.Cells.Item(268, 9).Value2 = "Primera"
$ws.Cells.Item(268, 10).Value2 = 700
$ws.Cells.Item(268, 11).Value2 = 17000
$ws.Cells.Item(268, 12).Value2 = 18000
$ws.Cells.Item(268, 13).Value2 = 17571
$ws.Cells.Item(268, 14).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(268, 15).Value2 = "Limache"
$ws.Cells.Item(268, 16).Value2 = 976
$ws.Cells.Item(268, 17).Value2 = 18
$ws.Cells.Item(268, 18).Value2 = "Hortaliza"

# --- New row 269 ---
$ws.Cells.Item(269, 1).Value2 = 11
$ws.Cells.Item(269, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(269, 3).Value2 = "Bíobío"
$ws.Cells.Item(269, 4).Value2 = 44509
$ws.Cells.Item(269, 5).Value2 = 8
$ws.Cells.Item(269, 6).Value2 = 100112020
$ws.Cells.Item(269, 7).Value2 = "Tomate"
$ws.Cells.Item(269, 8).Value2 = "Larga vida"
$ws.Cells.Item(269, 9).Value2 = "Segunda"
$ws.Cells.Item(269, 10).Value2 = 400
$ws.Cells.Item(269, 11).Value2 = 16000
$ws.Cells.Item(269, 12).Value2 = 16000
$ws.Cells.Item(269, 13).Value2 = 16000
$ws.Cells.Item(269, 14).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(269, 15).Value2 = "Limache"
$ws.Cells.Item(269, 16).Value2 = 889
$ws.Cells.Item(269, 17).Value2 = 18
$ws.Cells.Item(269, 18).Value2 = "Hortaliza"
